# Add homework row for 2020-03-12 (row 73) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 73

$ws.Cells.Item($row, 1).Value = 1583971200

# Force B73/C73 to be stored as text (not auto-converted to a date serial
# / number by Excel's type inference), then restore the default "Normal"
# style so no stray number-format style is left behind on the cell.
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = "2020-03-12"
$ws.Cells.Item($row, 2).Style = "Normal"

$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "0215"
$ws.Cells.Item($row, 3).Style = "Normal"

$ws.Cells.Item($row, 4).Value = "SLVEST"
$ws.Cells.Item($row, 5).Value = 0.875
$ws.Cells.Item($row, 6).Value = 0.885
$ws.Cells.Item($row, 7).Value = 0.8
$ws.Cells.Item($row, 8).Value = 0.8100000000000001
$ws.Cells.Item($row, 9).Value = 13119600
